$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Image Streaming" (sheet2.xml)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Image Streaming")

# Header rename: D2 "Time Taken for sending each frame of image"
#             -> "Average Time Taken for sending each frame of image"
$ws2.Range("D2").Value = "Average Time Taken for sending each frame of image"

# New "Comment" header in F2
$ws2.Range("F2").Value = "Comment"

# Existing data edits: row4 (14 ms -> 15 ms), row5 (1 sec -> 75ms)
$ws2.Range("D4").Value = "15 ms"
$ws2.Range("D5").Value = "75ms"

# New rows 6-8 (S.NO 4,5,6) for 1920x1080 resolution
$ws2.Range("B6").Value = 4
$ws2.Range("C6").Value = "1920x1080"
$ws2.Range("D6").Value = "10 ms"
$ws2.Range("E6").Value = "Compressed string to byte array"

$ws2.Range("B7").Value = 5
$ws2.Range("C7").Value = "1920x1080"
$ws2.Range("D7").Value = "30 ms"
$ws2.Range("E7").Value = "U32 1D array"

$ws2.Range("B8").Value = 6
$ws2.Range("C8").Value = "1920x1080"
$ws2.Range("D8").Value = "1.3 s"
$ws2.Range("E8").Value = "Json string"

# Match B-column style (left aligned, no border) used by the other S.NO cells
$ws2.Range("B6:B8").HorizontalAlignment = -4131

# New "Comment" column: merged F3:F8, centered, value "30 fps"
$ws2.Range("F3:F8").Merge()
$ws2.Range("F3").Value = "30 fps"
$ws2.Range("F3:F8").HorizontalAlignment = -4108

# Column F width + row 2 height
$ws2.Columns("F").ColumnWidth = 23.28515625
$ws2.Rows(2).RowHeight = 62.25

# Selection / active cell
$ws2.Activate()
$ws2.Range("K5").Select()

# ---------------------------------------------------------------------------
# Sheet "Image Transfer to IS" (sheet1.xml)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Image Transfer to IS")
$ws1.Activate()
$ws1.Range("I18").Select()
